# Generate Report for Handoff
# Adds two new handed-off files (4358dd1c-... and 591e6bd7-...) as new rows
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the layout
# of the existing rows (2-5) on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Matches the blue-underline "HyperLink" cell style (font) already used by
# column A on every sheet for linked file names.
function Set-HyperlinkFont($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# New item 1: 4358dd1c-fff6-4102-864b-305cb2428e60
# ---------------------------------------------------------------------
$uuid1 = "4358dd1c-fff6-4102-864b-305cb2428e60"
$rev1  = "c80f05e1b798a5591cbe3d02cce861727c176007"

# Overview sheet, row 6
$wsOverview.Range("A6").Value = "$uuid1.md"
Set-HyperlinkFont $wsOverview.Range("A6")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev1/e2e/$uuid1.md", "", "", "$uuid1.md")
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-28-17 12:28:35"

# zh-cn sheet, row 6
$wsZhCn.Range("A6").Value = "$uuid1.md"
Set-HyperlinkFont $wsZhCn.Range("A6")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev1/e2e/$uuid1.md", "", "", "$uuid1.md")
$wsZhCn.Range("B6").Value = ".md"
Set-HyperlinkFont $wsZhCn.Range("B6")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev1/e2e/$uuid1.md", "", "", ".md")
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "$uuid1.$rev1.zh-cn.xlf"
Set-HyperlinkFont $wsZhCn.Range("D6")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$rev1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$uuid1.$rev1.zh-cn.xlf", "", "", "$uuid1.$rev1.zh-cn.xlf")
$wsZhCn.Range("E6").Value = "2016-03-17 12:28:32"
$wsZhCn.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I6").Value = "Include"

# de-de sheet, row 6
$wsDeDe.Range("A6").Value = "$uuid1.md"
Set-HyperlinkFont $wsDeDe.Range("A6")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev1/e2e/$uuid1.md", "", "", "$uuid1.md")
$wsDeDe.Range("B6").Value = ".md"
Set-HyperlinkFont $wsDeDe.Range("B6")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev1/e2e/$uuid1.md", "", "", ".md")
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "$uuid1.$rev1.de-de.xlf"
Set-HyperlinkFont $wsDeDe.Range("D6")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$rev1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$uuid1.$rev1.de-de.xlf", "", "", "$uuid1.$rev1.de-de.xlf")
$wsDeDe.Range("E6").Value = "2016-03-17 12:28:35"
$wsDeDe.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I6").Value = "Include"

# ---------------------------------------------------------------------
# New item 2: 591e6bd7-243e-4d56-b3dd-c8d740216bf9
# ---------------------------------------------------------------------
$uuid2 = "591e6bd7-243e-4d56-b3dd-c8d740216bf9"
$rev2  = "a8ea106679136a2a56c20e85a48901152851559d"

# Overview sheet, row 7
$wsOverview.Range("A7").Value = "$uuid2.md"
Set-HyperlinkFont $wsOverview.Range("A7")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev2/e2e/$uuid2.md", "", "", "$uuid2.md")
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-28-17 12:28:35"

# zh-cn sheet, row 7
$wsZhCn.Range("A7").Value = "$uuid2.md"
Set-HyperlinkFont $wsZhCn.Range("A7")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev2/e2e/$uuid2.md", "", "", "$uuid2.md")
$wsZhCn.Range("B7").Value = ".md"
Set-HyperlinkFont $wsZhCn.Range("B7")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev2/e2e/$uuid2.md", "", "", ".md")
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "$uuid2.$rev2.zh-cn.xlf"
Set-HyperlinkFont $wsZhCn.Range("D7")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$rev2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$uuid2.$rev2.zh-cn.xlf", "", "", "$uuid2.$rev2.zh-cn.xlf")
$wsZhCn.Range("E7").Value = "2016-03-17 12:28:32"
$wsZhCn.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I7").Value = "Include"

# de-de sheet, row 7
$wsDeDe.Range("A7").Value = "$uuid2.md"
Set-HyperlinkFont $wsDeDe.Range("A7")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev2/e2e/$uuid2.md", "", "", "$uuid2.md")
$wsDeDe.Range("B7").Value = ".md"
Set-HyperlinkFont $wsDeDe.Range("B7")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/$rev2/e2e/$uuid2.md", "", "", ".md")
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "$uuid2.$rev2.de-de.xlf"
Set-HyperlinkFont $wsDeDe.Range("D7")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$rev2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$uuid2.$rev2.de-de.xlf", "", "", "$uuid2.$rev2.de-de.xlf")
$wsDeDe.Range("E7").Value = "2016-03-17 12:28:35"
$wsDeDe.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I7").Value = "Include"
